# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Austria (row 35) ---
$ws.Range("B35").Value = 15882
$ws.Range("C35").Value = 11
$ws.Range("D35").Value = 14061
$ws.Range("E35").Value = 1201
$ws.Range("F35").Value = 68
$ws.Range("G35").Value = 2
$ws.Range("H35").Value = 620

# --- Marruecos (row 55) ---
$ws.Range("B55").Value = 6226
$ws.Range("C55").Value = 163
$ws.Range("D55").Value = 2759
$ws.Range("E55").Value = 3279

# --- Albania (row 103) ---
$ws.Range("B103").Value = 872
$ws.Range("C103").Value = 4
$ws.Range("D103").Value = 654

# --- San Marino (row 119) ---
$ws.Range("D119").Value = 130
$ws.Range("E119").Value = 457
$ws.Range("F119").Value = 2

# --- Belice / Nueva Caledonia swap positions (rows 192-193) ---
$ws.Range("A192").Value = "Nueva Caledonia"
$ws.Range("D192").Value = 18
$ws.Range("H192").Value = 0

$ws.Range("A193").Value = "Belice"
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2
